$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Columns("D:D").Insert()

# Copy number formats/styles from the (now shifted) E column into the new D column
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D with the new quarter (2018-09-30) figures
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 148800
$ws.Range("D9").Value = 74900
$ws.Range("D10").Value = 73900
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -3700
$ws.Range("D15").Value = 24600
$ws.Range("D17").Value = 111600
$ws.Range("D18").Value = 37200
$ws.Range("D20").Value = -4200
$ws.Range("D21").Value = 33100
$ws.Range("D22").Value = 32700
$ws.Range("D23").Value = 200
$ws.Range("D24").Value = -400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 600
$ws.Range("D27").Value = -4400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 4200
$ws.Range("D33").Value = -4400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -4400
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 132900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 69800
$ws.Range("D44").Value = 29300
$ws.Range("D45").Value = 44800
$ws.Range("D46").Value = 276700
$ws.Range("D47").Value = 201500
$ws.Range("D48").Value = 1695500
$ws.Range("D49").Value = 272600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 42500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2488900
$ws.Range("D57").Value = 78700
$ws.Range("D58").Value = 32900
$ws.Range("D59").Value = 121400
$ws.Range("D60").Value = 232900
$ws.Range("D61").Value = 1566500
$ws.Range("D62").Value = 133200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2037200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -234000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 451600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -4400
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 10100
$ws.Range("D91").Value = -2400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 74600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -63100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 21500